# Actualización 11 de Mayo - Mañana
$wb = $excel.ActiveWorkbook

# --- Sheet "Estadisticos 2P": update second-partial statistics row ---
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws2.Range("D2").Value = 16
$ws2.Range("E2").Value = 6
$ws2.Range("F2").Value = 23
$ws2.Range("G2").Value = 58.97

# --- Sheet "Rescatables": replace the two retake rows and add a third ---
$ws4 = $wb.Worksheets.Item("Rescatables")

# Row 2 (was GARCIA / GUZMAN / ALEXIS MANUEL, NC 20330051920014)
$ws4.Range("A2").Value = 20330051920013
$ws4.Range("B2").Value = "FLORES"
$ws4.Range("C2").Value = "JUAREZ"
$ws4.Range("D2").Value = "LUIS ANGEL"
$ws4.Range("E2").Value = "REALIZA INSTALACIONES ELÉCTRICAS"
$ws4.Range("F2").Value = "2AEM"
$ws4.Range("G2").Value = 2

# Row 3 (was HERNANDEZ / BALDERAS / JUAN MANUEL, NC 20330051920015)
$ws4.Range("A3").Value = 20330051920030
$ws4.Range("B3").Value = "TZANAHUA"
$ws4.Range("C3").Value = "GONZALEZ"
$ws4.Range("D3").Value = "XIMENA"
$ws4.Range("E3").Value = "REALIZA INSTALACIONES ELÉCTRICAS"
$ws4.Range("F3").Value = "2AEM"
$ws4.Range("G3").Value = 2

# Row 4 (new)
$ws4.Range("A4").Value = 20330051920014
$ws4.Range("B4").Value = "GARCIA"
$ws4.Range("C4").Value = "GUZMAN"
$ws4.Range("D4").Value = "ALEXIS MANUEL"
$ws4.Range("E4").Value = "REALIZA INSTALACIONES ELÉCTRICAS"
$ws4.Range("F4").Value = "2AEM"
$ws4.Range("G4").Value = 1
